# Update the "dSF" (column F) values for specific rows as per the
# repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -8
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = -4
$ws.Range("F10").Value = -3
$ws.Range("F14").Value = 10
$ws.Range("F19").Value = -10
